# Fill in the "Why important (one sentence)" notes for the Models@run.time
# subchapter table (sheet "Models@runtime"), rows 3-6, column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$e3 = @"
Models at runtime aim to facilitate the process of dynamically adding new requirements to software by adding an extra reflection level that maps a model to the underlying software. As the author states, models at runtime make use of MDE principles.
Models at runtime aim to blur the distinction between software design, its evolution, and execution.
"@

$e4 = @"
Runtime modeling connects software with its abstraction to bidirectionally synchronize an artifact with its reflection. Consequently, it lets to reduce the scheduled downtime of a system since semantic changes and behavior modification happen at runtime. 
Continuing the main idea of MDE to view models as the most important artifacts during development, models@runtime encourages the use of models not only to design software but also during its execution. Therefore, the running system and its model evolve concurrently. 
An example that showcases one of the advantages is described below. Consider a system the domain diagram of which is present below. It is not uncommon that after some time, a new requirement is added. In a classical scenario, a development team would have to stop the software, adapt the domain model, regenerate code if needed and deploy the system again. All those cumbersome stages can be avoided by using models at runtime to adapt the domain model on-the-fly. In can be achieved by running models using an execution engine that reflects the current state of a system.
"@

$e5 = @"
There is an increasing demand for self-adaptive systems that can deal with unexpected software changes. As stated here, the three main pillars of models at runtime are modeling, separation of concerns, and reflection.
The architecture of a model at runtime system is depicted in figure 1. According to the source, it consists of three interrelated parts.
The first level entitles the models of a target system. It has four model subtypes that abstract the target system. Context models contain the current state of the system environment. Configuration models depict the architectural outlook on the running system. Capability models contain attributes to manage the target system. Even though this model is mostly static and depends on the running infrastructure, it still can be updated after a new component enters a system. The last type of model is a planning model that contains instructions on embedding new components into a target system.
The next level serves configurational purposes and includes system components that enable feedback interconnection between this management level and the level below. It has active entities such as a reasoner and analyzer to manipulate the models defined on a base level. The reasoner reasons about the configuration state of the system in the future. The analyzer checks whether the target system state corresponds to the goals of the system. If it is not the case, then the analyzer fires a reasoner to reevaluate the configuration. Optionally, a learner can be defined on this level to maintain the synchronization between the models of the base level with the target system and check the usefulness of the reasoner's statements.
The top level comprises goal-related models of a system and enables connection with the configuration management layer. Those types of models are fed by the reasoner to check if the goals can be fulfilled. It is worth mentioning, that such models tend to evolve over time in the presence of new requirements and goals.
"@

$e6 = @"
The classification of used models at runtime contains three categories. Namely, models can be categorized by the purpose of their use, the language by which they are constructed, and, finally, by the area they represent.
"@

# Trim the single trailing newline the here-strings pick up before "@
$e3 = $e3.TrimEnd("`r", "`n")
$e4 = $e4.TrimEnd("`r", "`n")
$e5 = $e5.TrimEnd("`r", "`n")
$e6 = $e6.TrimEnd("`r", "`n")

$ws.Range("E3").Value = $e3
$ws.Range("E4").Value = $e4
$ws.Range("E5").Value = $e5
$ws.Range("E6").Value = $e6

# Update the saved selection to match the author's final cursor position.
$ws.Activate()
$ws.Range("E4").Select()
